# Updated cryptos list on Wed Jun 12 10:30:35 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.850.86'
$ws.Range("E2").Value = '  +0.89%  '

$ws.Range("D3").Value = '3.546.09'
$ws.Range("E3").Value = '  +0.05%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.29%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '616.60'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.76%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '152.65'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.53%  '

$ws.Range("D7").Value = '3.543.23'
$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("E8").Value = '  +0.11%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.482'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.90%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.140'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.23%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.08'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.51%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.427'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.10%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000221'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.27%  '

$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '4.152.81'
$ws.Range("E14").Value = '  +0.41%  '

$ws.Range("B15").Value = 'Avalanche'
$ws.Range("C15").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '32.21'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.01%  '

$ws.Range("D16").Value = '3.551.48'
$ws.Range("E16").Value = '  +0.69%  '

$ws.Range("D17").Value = '67.668.59'
$ws.Range("E17").Value = '  +0.88%  '

$ws.Range("E18").Value = '  -0.61%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.42'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.17%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.39'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.86%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.71'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.08%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '448.05'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.40%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.626'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.04%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '77.53'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.97%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000132'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +5.89%  '

$ws.Range("D26").Value = '3.690.14'
$ws.Range("E26").Value = '  +0.56%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.33'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.34%  '

$ws.Range("E28").Value = '  +0.06%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.66'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.11%  '

$ws.Range("E30").Value = '  -1.39%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.60'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.45%  '

$ws.Range("E32").Value = '  +6.07%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.20%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.01'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.18%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.22'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.21%  '

$ws.Range("D36").Value = '3.539.37'
$ws.Range("E36").Value = '  +0.55%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.86'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.43%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.07'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.30%  '

$ws.Range("B40").Value = 'FirstDigitalUSD'
$ws.Range("C40").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.47%  '

$ws.Range("B41").Value = 'Monero'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '176.97'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.39%  '

$ws.Range("B42").Value = 'Hedera'
$ws.Range("C42").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0898'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.82%  '

$ws.Range("B43").Value = 'Stacks'
$ws.Range("C43").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.18'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.47%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.44'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.59%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.887'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.89%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '28.66'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.68%  '

$ws.Range("E47").Value = '  -0.77%  '

$ws.Range("E48").Value = '  -1.96%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.28'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.92%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.63'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.83%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.999'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.25%  '
